$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 581, shifting existing rows 581-690 down to 582-691
$ws.Rows.Item(581).Insert()

# Populate the new row 581 with data
$ws.Cells.Item(581, 1).Value = 11
$ws.Cells.Item(581, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(581, 3).Value = "Bíobío"
$ws.Cells.Item(581, 4).Value = 45258
$ws.Cells.Item(581, 5).Value = 8
$ws.Cells.Item(581, 6).Value = 100112006
$ws.Cells.Item(581, 7).Value = "Repollo"
$ws.Cells.Item(581, 8).Value = "Crespo record"
$ws.Cells.Item(581, 9).Value = "Primera"
$ws.Cells.Item(581, 10).Value = 1000
$ws.Cells.Item(581, 11).Value = 1000
$ws.Cells.Item(581, 12).Value = 1000
$ws.Cells.Item(581, 13).Value = 1000
$ws.Cells.Item(581, 14).Value = "$/unidad"
$ws.Cells.Item(581, 15).Value = "Región Metropolitana"
$ws.Cells.Item(581, 16).Value = 1000
$ws.Cells.Item(581, 17).Value = 1
$ws.Cells.Item(581, 18).Value = "Hortaliza"

# Copy the style of column D (date-formatted) from the row above into the new row's D cell
$ws.Cells.Item(580, 4).Copy()
$ws.Cells.Item(581, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(581, 4).Value = 45258
